$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("PayNowCC")
$ws.Range("B2").Value = "Fri Nov 29 18:09:38 IST 2024"
$ws.Range("B3").Value = "Fri Nov 29 18:27:04 IST 2024"
$ws.Range("B4").Value = "Fri Nov 29 18:44:34 IST 2024"
$ws.Range("B5").Value = "Fri Nov 29 19:02:10 IST 2024"
$ws.Range("B6").Value = "Fri Nov 29 19:19:18 IST 2024"

$ws = $wb.Worksheets.Item("PayNowCorpSCF")
$ws.Range("A2").Value = "Fail"
$ws.Range("B2").Value = "Fri Nov 29 09:33:45 IST 2024"
$ws.Range("A3").Value = "Fail"
$ws.Range("B3").Value = "Fri Nov 29 09:54:09 IST 2024"
$ws.Range("A4").Value = "Fail"
$ws.Range("B4").Value = "Fri Nov 29 10:13:44 IST 2024"
$ws.Range("A5").Value = "Fail"
$ws.Range("B5").Value = "Fri Nov 29 10:32:22 IST 2024"

$ws = $wb.Worksheets.Item("PayNowCreditSCF")
$ws.Range("A2").Value = "Fail"
$ws.Range("B2").Value = "Fri Nov 29 12:04:58 IST 2024"
$ws.Range("A3").Value = "Fail"
$ws.Range("B3").Value = "Fri Nov 29 12:24:35 IST 2024"
$ws.Range("B4").Value = "Fri Nov 29 12:44:00 IST 2024"
$ws.Range("B5").Value = "Fri Nov 29 13:01:20 IST 2024"

$ws = $wb.Worksheets.Item("PayNowPersonalCheckSCF")
$ws.Range("B2").Value = "Fri Nov 29 14:32:03 IST 2024"
$ws.Range("B3").Value = "Fri Nov 29 14:49:14 IST 2024"
$ws.Range("A4").Value = "Fail"
$ws.Range("B4").Value = "Fri Nov 29 15:07:02 IST 2024"
$ws.Range("A5").Value = "Fail"
$ws.Range("B5").Value = "Fri Nov 29 15:25:51 IST 2024"

$ws = $wb.Worksheets.Item("PayNowPersonalSavingsSCF")
$ws.Range("B2").Value = "Fri Nov 29 16:56:53 IST 2024"
$ws.Range("B3").Value = "Fri Nov 29 17:14:06 IST 2024"
$ws.Range("A4").Value = "Fail"
$ws.Range("B4").Value = "Fri Nov 29 17:32:08 IST 2024"
$ws.Range("A5").Value = "Fail"
$ws.Range("B5").Value = "Fri Nov 29 17:50:39 IST 2024"

$ws = $wb.Worksheets.Item("PayNowCreditDCF")
$ws.Range("A2").Value = "Fail"
$ws.Range("B2").Value = "Fri Nov 29 10:50:45 IST 2024"
$ws.Range("A3").Value = "Fail"
$ws.Range("B3").Value = "Fri Nov 29 11:10:26 IST 2024"
$ws.Range("B4").Value = "Fri Nov 29 11:30:10 IST 2024"
$ws.Range("B5").Value = "Fri Nov 29 11:47:34 IST 2024"

$ws = $wb.Worksheets.Item("PayNowCorpDCF")
$ws.Range("A2").Value = "Fail"
$ws.Range("B2").Value = "Fri Nov 29 06:27:41 IST 2024"
$ws.Range("A3").Value = "Fail"
$ws.Range("B3").Value = "Fri Nov 29 06:47:46 IST 2024"
$ws.Range("A4").Value = "Fail"
$ws.Range("B4").Value = "Fri Nov 29 07:06:59 IST 2024"
$ws.Range("A5").Value = "Fail"
$ws.Range("B5").Value = "Fri Nov 29 07:25:16 IST 2024"

$ws = $wb.Worksheets.Item("PayNowPersonalCheckDCF")
$ws.Range("B2").Value = "Fri Nov 29 15:44:26 IST 2024"
$ws.Range("B3").Value = "Fri Nov 29 16:01:51 IST 2024"
$ws.Range("A4").Value = "Fail"
$ws.Range("B4").Value = "Fri Nov 29 16:18:54 IST 2024"
$ws.Range("A5").Value = "Fail"
$ws.Range("B5").Value = "Fri Nov 29 16:37:45 IST 2024"

$ws = $wb.Worksheets.Item("OverAndUnderPayCredit")
$ws.Range("A2").Value = "Fail"
$ws.Range("B2").Value = "Fri Nov 29 02:50:41 IST 2024"
$ws.Range("A3").Value = "Fail"
$ws.Range("B3").Value = "Fri Nov 29 03:09:19 IST 2024"
$ws.Range("B4").Value = "Fri Nov 29 03:27:31 IST 2024"
$ws.Range("B5").Value = "Fri Nov 29 03:44:59 IST 2024"

$ws = $wb.Worksheets.Item("OverAndUnderPayPC")
$ws.Range("A2").Value = "Fail"
$ws.Range("B2").Value = "Fri Nov 29 04:01:58 IST 2024"
$ws.Range("A3").Value = "Fail"
$ws.Range("B3").Value = "Fri Nov 29 04:20:29 IST 2024"
$ws.Range("A4").Value = "Fail"
$ws.Range("B4").Value = "Fri Nov 29 04:38:49 IST 2024"
$ws.Range("A5").Value = "Fail"
$ws.Range("B5").Value = "Fri Nov 29 04:56:20 IST 2024"

$ws = $wb.Worksheets.Item("OverAndUnderPayPS")
$ws.Range("A2").Value = "Fail"
$ws.Range("B2").Value = "Fri Nov 29 05:14:24 IST 2024"
$ws.Range("A3").Value = "Fail"
$ws.Range("B3").Value = "Fri Nov 29 05:32:59 IST 2024"
$ws.Range("A4").Value = "Fail"
$ws.Range("B4").Value = "Fri Nov 29 05:51:36 IST 2024"
$ws.Range("A5").Value = "Fail"
$ws.Range("B5").Value = "Fri Nov 29 06:09:39 IST 2024"

$ws = $wb.Worksheets.Item("OverAndUnderPayCorp")
$ws.Range("A2").Value = "Fail"
$ws.Range("B2").Value = "Fri Nov 29 01:38:50 IST 2024"
$ws.Range("A3").Value = "Fail"
$ws.Range("B3").Value = "Fri Nov 29 01:57:27 IST 2024"
$ws.Range("A4").Value = "Fail"
$ws.Range("B4").Value = "Fri Nov 29 02:15:35 IST 2024"
$ws.Range("A5").Value = "Fail"
$ws.Range("B5").Value = "Fri Nov 29 02:33:03 IST 2024"

$ws = $wb.Worksheets.Item("NoModifyAmountCorp")
$ws.Range("A2").Value = "Fail"
$ws.Range("B2").Value = "Thu Nov 28 23:53:07 IST 2024"
$ws.Range("A3").Value = "Fail"
$ws.Range("B3").Value = "Fri Nov 29 00:11:00 IST 2024"

$ws = $wb.Worksheets.Item("NoModifyAmountPC")
$ws.Range("A2").Value = "Fail"
$ws.Range("B2").Value = "Fri Nov 29 00:28:23 IST 2024"
$ws.Range("A3").Value = "Fail"
$ws.Range("B3").Value = "Fri Nov 29 00:46:16 IST 2024"

$ws = $wb.Worksheets.Item("NoModifyAmountPS")
$ws.Range("A2").Value = "Fail"
$ws.Range("B2").Value = "Fri Nov 29 01:03:30 IST 2024"
$ws.Range("A3").Value = "Fail"
$ws.Range("B3").Value = "Fri Nov 29 01:21:29 IST 2024"

$ws = $wb.Worksheets.Item("NoModifyAmountCC")
$ws.Range("A2").Value = "Fail"
$ws.Range("B2").Value = "Fri Nov 29 21:35:12 IST 2024"
$ws.Range("B3").Value = "Thu Nov 28 23:36:16 IST 2024"

$ws = $wb.Worksheets.Item("PayNowCorp")
$ws.Range("A2").Value = "Fail"
$ws.Range("B2").Value = "Fri Nov 29 07:44:00 IST 2024"
$ws.Range("A3").Value = "Fail"
$ws.Range("B3").Value = "Fri Nov 29 08:02:15 IST 2024"
$ws.Range("A4").Value = "Fail"
$ws.Range("B4").Value = "Fri Nov 29 08:21:24 IST 2024"
$ws.Range("A5").Value = "Fail"
$ws.Range("B5").Value = "Fri Nov 29 08:39:42 IST 2024"
$ws.Range("A6").Value = "Fail"
$ws.Range("B6").Value = "Fri Nov 29 08:57:12 IST 2024"
$ws.Range("A7").Value = "Fail"
$ws.Range("B7").Value = "Fri Nov 29 09:15:25 IST 2024"
